$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text, matching the source data which
# stores price strings like "1.000" / "27.877.37" as inline text, not numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$priceUpdates = @{
    2 = "27.877.37"
    3 = "1.908.30"
    4 = "1.000"
    5 = "313.00"
    6 = "1.000"
    7 = "0.5012"
    8 = "0.3806"
    9 = "0.07293"
    10 = "0.9090"
    11 = "21.22"
    12 = "0.07654"
    13 = "1.899.00"
    14 = "5.478"
    15 = "92.54"
    16 = "1.001"
    17 = "0.000008731"
    18 = "0.9993"
    19 = "27.895.32"
    20 = "14.64"
    21 = "5.171"
    22 = "2.126.12"
    23 = "10.85"
    24 = "6.613"
    25 = "153.02"
    26 = "1.839"
    27 = "2.198"
    28 = "18.40"
    29 = "115.25"
    30 = "4.909"
    31 = "0.09017"
    32 = "3.202"
    33 = "4.823"
    34 = "1.232"
    35 = "0.7783"
    36 = "2.621"
    37 = "0.02081"
    38 = "3.076"
    39 = "1.092"
    40 = "0.5542"
    41 = "0.05265"
    42 = "6.814"
    43 = "114.24"
    44 = "8.509"
    45 = "0.1517"
    46 = "0.4813"
    47 = "10.55"
    48 = "1.000"
    49 = "1.634"
    50 = "67.30"
    51 = "0.06051"
}

$volumeUpdates = @{
    2 = "  -0.51%  "
    3 = "  +0.00%  "
    4 = "  -0.55%  "
    5 = "  -1.62%  "
    6 = "  -0.49%  "
    7 = "  +3.76%  "
    8 = "  +0.00%  "
    9 = "  -0.91%  "
    10 = "  -2.71%  "
    11 = "  +1.91%  "
    12 = "  -1.78%  "
    13 = "  -0.55%  "
    14 = "  -0.37%  "
    15 = "  +0.72%  "
    16 = "  -0.55%  "
    17 = "  -1.73%  "
    18 = "  -0.51%  "
    19 = "  -0.55%  "
    20 = "  -0.60%  "
    21 = "  +0.25%  "
    22 = "  -1.71%  "
    23 = "  -0.51%  "
    24 = "  -0.45%  "
    25 = "  -2.68%  "
    26 = "  -3.96%  "
    27 = "  +3.55%  "
    28 = "  -0.60%  "
    29 = "  -1.58%  "
    30 = "  -1.56%  "
    31 = "  +0.78%  "
    32 = "  -2.57%  "
    33 = "  +3.39%  "
    34 = "  -1.86%  "
    35 = "  +0.74%  "
    36 = "  +0.66%  "
    37 = "  +1.26%  "
    38 = "  +2.69%  "
    39 = "  -1.52%  "
    40 = "  +0.31%  "
    41 = "  -0.54%  "
    42 = "  -2.64%  "
    43 = "  +3.61%  "
    44 = "  +0.08%  "
    45 = "  -0.81%  "
    46 = "  -0.27%  "
    47 = "  -1.04%  "
    48 = "  -0.51%  "
    49 = "  -0.86%  "
    50 = "  -1.26%  "
    51 = "  -0.45%  "
}

foreach ($row in $priceUpdates.Keys) {
    $ws.Cells.Item($row, 4).Value = $priceUpdates[$row]
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Cells.Item($row, 5).Value = $volumeUpdates[$row]
}
